# Regenerate the experiment task-order sheets: new randomized stimulus file
# lists per task, and the sheets shuffled into a new tab order.
#
# The workbook always has five task-order sheets (GNG, NB, RS, TOL, vSAT).
# Each keeps its own identity/content-shape (row count never changes for a
# given task), but gets a freshly generated name suffix, freshly generated
# stimulus-file values in column B, and is moved to a new tab position.

$wb = $excel.ActiveWorkbook

# Grab the five original sheets by their current (original) tab order.
$sGNG  = $wb.Worksheets.Item(1)
$sNB   = $wb.Worksheets.Item(2)
$sRS   = $wb.Worksheets.Item(3)
$sTOL  = $wb.Worksheets.Item(4)
$sVSAT = $wb.Worksheets.Item(5)

# --- Rename each to its freshly generated sheet name -----------------------
$sGNG.Name  = "GNG_TO-16515889306453595"
$sNB.Name   = "NB_TO-16515889304667995"
$sRS.Name   = "RS_TO-16515889305982409"
$sTOL.Name  = "TOL_TO-1651588930515023"
$sVSAT.Name = "vSAT_TO-1651588930595209"

# --- Reorder tabs to the new order: NB, TOL, vSAT, RS, GNG ------------------
# (Worksheets.Item(N) handles rebind to whatever occupies that slot after a
# move, so each sheet is re-looked-up by its now-stable name right before
# it's moved. Walking the target order back-to-front and always moving the
# next sheet to the very front reproduces the desired final order.)
$wb.Worksheets.Item("GNG_TO-16515889306453595").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("RS_TO-16515889305982409").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("vSAT_TO-1651588930595209").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("TOL_TO-1651588930515023").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("NB_TO-16515889304667995").Move($wb.Worksheets.Item(1))

# --- NB sheet: 9 stimulus rows ----------------------------------------------
$ws = $wb.Worksheets.Item("NB_TO-16515889304667995")
$ws.Range("B2").Value = "ZB-match_3-16515889292623913.csv"
$ws.Range("B3").Value = "OB-16515889294041233.csv"
$ws.Range("B4").Value = "OB-16515889296314902.csv"
$ws.Range("B5").Value = "ZB-match_9-16515889287934875.csv"
$ws.Range("B6").Value = "ZB-match_4-16515889293092864.csv"
$ws.Range("B7").Value = "TB-1651588929938002.csv"
$ws.Range("B8").Value = "OB-16515889295210896.csv"
$ws.Range("B9").Value = "TB-16515889304361045.csv"
$ws.Range("B10").Value = "TB-16515889300218463.csv"

# --- TOL sheet: 6 stimulus rows ---------------------------------------------
$ws = $wb.Worksheets.Item("TOL_TO-1651588930515023")
$ws.Range("B2").Value = "MM_stims-16515889304826152.csv"
$ws.Range("B3").Value = "ZM_stims-16515889304707072.csv"
$ws.Range("B4").Value = "MM_stims-16515889304980178.csv"
$ws.Range("B5").Value = "ZM_stims-16515889304836328.csv"
$ws.Range("B6").Value = "MM_stims-16515889305139832.csv"
$ws.Range("B7").Value = "ZM_stims-16515889304990134.csv"

# --- vSAT sheet: 4 stimulus rows --------------------------------------------
$ws = $wb.Worksheets.Item("vSAT_TO-1651588930595209")
$ws.Range("B2").Value = "SAT_stims-1651588930522113.csv"
$ws.Range("B3").Value = "vSAT_stims-16515889305798962.csv"
$ws.Range("B4").Value = "vSAT_stims-16515889305625527.csv"
$ws.Range("B5").Value = "SAT_stims-16515889305471237.csv"

# --- RS sheet: unchanged content (eyes closed / eyes open) -----------------
# (no cell-value changes needed, only the rename/move above)

# --- GNG sheet: 4 stimulus rows ---------------------------------------------
$ws = $wb.Worksheets.Item("GNG_TO-16515889306453595")
$ws.Range("B2").Value = "go_stims-1651588930601185.csv"
$ws.Range("B3").Value = "GNG_stims-16515889306280458.csv"
$ws.Range("B4").Value = "go_stims-16515889306300511.csv"
$ws.Range("B5").Value = "GNG_stims-16515889306433635.csv"
